# The document ends with a "Requisitos" section followed by a footer-like
# block that was removed from the site build:
#   (empty paragraph)
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages.
#    Original theme under Creative Commons Attribution"
# These three paragraphs are deleted, leaving the trailing empty paragraph
# and the page-break paragraph untouched.

$d = $word.ActiveDocument

# Locate the requirement paragraph that anchors the block to remove.
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("LOB1039: Física Experimental III (Requisito fraco)",
                      $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
if (-not $anchor.Find.Found) {
    throw "Anchor paragraph not found"
}

# Locate the copyright/footer paragraph that ends the block to remove.
$footer = $d.Content.Duplicate
$footer.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
                      $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
if (-not $footer.Find.Found) {
    throw "Footer paragraph not found"
}

# Resolve the actual Paragraph objects/indices owning those ranges by
# scanning the document's Paragraphs collection (Range.Paragraphs.First /
# .Next are not reliable on duplicated Find ranges in this runtime).
$anchorIndex = -1
$footerIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Start -eq $anchor.Start) {
        $anchorIndex = $i
    }
    if ($p.Range.Start -eq $footer.Start) {
        $footerIndex = $i
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not resolve anchor paragraph index"
}
if ($footerIndex -eq -1) {
    throw "Could not resolve footer paragraph index"
}

# The paragraph right after the anchor starts the block being removed
# (the blank separator paragraph); the footer paragraph ends it.
$startPara = $d.Paragraphs.Item($anchorIndex + 1)
$endPara = $d.Paragraphs.Item($footerIndex)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
